$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # RUNMANAGER
$ws3 = $wb.Worksheets.Item(3)   # USERACCOUNTMANAGEMENTDATA

# --- RUNMANAGER: add new test case row (row 9) ---
$ws1.Range("A9").Value = "verifyThatUserCannotChangePasswordWhenThePasswoedDoesNotMeetSpecifiedCriteria"
$ws1.Range("B9").Value = "To check this test is executed"
$ws1.Range("C9").Value = "yes"
$ws1.Range("D9").Value = "'8"
$ws1.Range("E9").Value = "'1"

# --- USERACCOUNTMANAGEMENTDATA: add new test data rows (rows 8 and 9) ---
$ws3.Range("A8").Value = "verifyThatUserCannotChangePasswordWhenThePasswoedDoesNotMeetSpecifiedCriteria"
$ws3.Range("B8").Value = "yes"
$ws3.Range("C8").Value = "Admin"
$ws3.Range("D8").Value = "admin123"
$ws3.Range("E8").Value = "Sunil"
$ws3.Range("F8").Value = "chrome"
$ws3.Range("G8").Value = "admin123"
$ws3.Range("H8").Value = "admin1"
$ws3.Range("I8").Value = "admin132"

$ws3.Range("A9").Value = "verifyThatUserCannotChangePasswordWhenThePasswoedDoesNotMeetSpecifiedCriteria"
$ws3.Range("B9").Value = "yes"
$ws3.Range("C9").Value = "Admin"
$ws3.Range("D9").Value = "admin123"
$ws3.Range("E9").Value = "Sunil"
$ws3.Range("F9").Value = "firefox"
$ws3.Range("G9").Value = "admin123"
$ws3.Range("H9").Value = "admin1"
$ws3.Range("I9").Value = "admin132"

# --- Column A autofit on both sheets (content grew) ---
$ws1.Columns.Item(1).AutoFit()
$ws3.Columns.Item(1).AutoFit()

# --- Selection / active sheet / active tab state ---
$ws1.Range("A9").Select() | Out-Null
$ws3.Activate() | Out-Null
$ws3.Range("E11").Select() | Out-Null
